# Weekly fruit/vegetable price update:
# A new price record (week of 2021-10-05) is inserted at row 149, pushing the
# existing history (rows 149-168) down by one row to 150-169.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 149, shifting rows 149:168 -> 150:169
$ws.Rows(149).Insert()

# Populate the newly inserted row 149 with the new weekly data point.
$ws.Cells.Item(149, 1).Value = 5
$ws.Cells.Item(149, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(149, 3).Value = "Maule"
$ws.Cells.Item(149, 4).Value = 44474
$ws.Cells.Item(149, 5).Value = 7
$ws.Cells.Item(149, 6).Value = 100112003
$ws.Cells.Item(149, 7).Value = "Ajo"
$ws.Cells.Item(149, 8).Value = "Chino"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 300
$ws.Cells.Item(149, 11).Value = 16000
$ws.Cells.Item(149, 12).Value = 16000
$ws.Cells.Item(149, 13).Value = 16000
$ws.Cells.Item(149, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(149, 15).Value = "China"
$ws.Cells.Item(149, 16).Value = 1600
$ws.Cells.Item(149, 17).Value = 10
$ws.Cells.Item(149, 18).Value = "Hortaliza"
